$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row
$ws2.Range("A1").Value = "Purchase Order Number"
$ws2.Range("B1").Value = "Date"
$ws2.Range("C1").Value = "Vendor Name"
$ws2.Range("D1").Value = "Vendor Address"
$ws2.Range("E1").Value = "Shipping Name"
$ws2.Range("F1").Value = "Shipping Address"
$ws2.Range("G1").Value = "Total Amount"
$ws2.Range("H1").Value = "Delivery by Date"
$ws2.Range("I1").Value = "Payment Terms"
$ws2.Range("J1").Value = "Client Name"
$ws2.Range("K1").Value = "Client Address"

# Data row
$ws2.Range("A2").Value = "#99674"
$ws2.Range("B2").Value = "31/08/2023"
$ws2.Range("C2").Value = "Leoba Connections (Pty) Ltd"
$ws2.Range("D2").Value = "103 Monale Street Munsieville"
$ws2.Range("E2").Value = "Lindiwe Mahlangu"
$ws2.Range("F2").Value = "1748 S K Matseke Avenue, Munsieville, Krugersdorp, Johannesburg, 1739"

$ws2.Range("G2").Value = 398055.02
$ws2.Range("G2").NumberFormat = "#,##0.00"

# Delivery by Date: reuse the same date style/format as Sheet1!B2 (numFmtId 14)
$ws2.Range("H2").Value = 45269
$ws1.Range("B2").Copy() | Out-Null
$ws2.Range("H2").PasteSpecial(-4122) | Out-Null

$ws2.Range("I2").Value = "No deliveries will be accepted after 1500 without prior arrangement. 2. Site will not be held responsible for offloading any materials. 3. The rates quantities of this order are fixed. 4. Order discrepancies must be challenged before any deliveries start. 5. Motheo reserves the right to cancel or amend this order at any time."
$ws2.Range("J2").Value = "Muhammed Saley"
$ws2.Range("K2").Value = "1748 S K Matseke Avenue, Munsieville, Krugersdorp, Johannesburg, 1739"
